$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.325.09'
$ws.Range("E2").Value = '  -0.83%  '
$ws.Range("D3").Value = '1.651.86'
$ws.Range("E3").Value = '  -0.51%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.17'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.67%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.510'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.73%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.46'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.76%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.258'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0614'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.90%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0889'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.49%  '
$ws.Range("D12").Value = '1.885.21'
$ws.Range("E12").Value = '  -0.44%  '
$ws.Range("D13").Value = '1.646.22'
$ws.Range("E13").Value = '  -0.78%  '
$ws.Range("E14").Value = '  -1.17%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.568'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +3.55%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.46'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.47%  '
$ws.Range("D17").Value = '27.335.19'
$ws.Range("E17").Value = '  -0.74%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '230.77'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -6.16%  '
$ws.Range("E19").Value = '  -0.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.41'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.81%  '
$ws.Range("E21").Value = '  -0.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.36'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.55%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.44'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.35%  '
$ws.Range("E24").Value = '  +0.72%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.27'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.91%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.10'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.81'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.45%  '
$ws.Range("E28").Value = '  -0.10%  '
$ws.Range("E29").Value = '  +0.15%  '
$ws.Range("E31").Value = '  -3.93%  '
$ws.Range("E32").Value = '  -1.29%  '
$ws.Range("D33").Value = '1.430.49'
$ws.Range("E33").Value = '  -0.52%  '
$ws.Range("E34").Value = '  +0.06%  '
$ws.Range("E35").Value = '  +1.14%  '
$ws.Range("E36").Value = '  -0.46%  '
$ws.Range("E37").Value = '  -2.57%  '
$ws.Range("E38").Value = '  -1.30%  '
$ws.Range("E39").Value = '  -0.28%  '
$ws.Range("E40").Value = '  +0.72%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.56'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +2.79%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '65.07'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -5.82%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.22'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.89%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.787'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.49%  '
$ws.Range("D46").Value = '1.794.13'
$ws.Range("E46").Value = '  -0.36%  '
$ws.Range("E47").Value = '  -2.43%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.91'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.91%  '
$ws.Range("E49").Value = '  -1.46%  '
$ws.Range("E50").Value = '  +0.06%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.73'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.20%  '
